$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correction: A3 was stored as text "000002"; fix it to the numeric value 2 ---
$ws.Range("A3").Value = 2

# --- Add the missing row 4 (retrait des donnees non conformes -> re-add corrected row) ---
# A4 must stay text "000003" (leading zeros). Force text formatting just long
# enough to make the assignment stick as a string, then clear the formatting
# again so no extra style is left behind on the cell (matches target which
# carries no explicit style index on A4).
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "000003"
$ws.Range("A4").ClearFormats()

$ws.Range("B4").Value = "C00009"
$ws.Range("C4").Value = 10
